# Refresh cryptocurrency price/volume data (automated update)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "29.458.67"
$ws.Range("E2").Value = "  +1.13%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.840.08"
$ws.Range("E3").Value = "  +0.28%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.001"
$ws.Range("E4").Value = "  -0.20%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "243.66"
$ws.Range("E5").Value = "  +0.30%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.6269"
$ws.Range("E6").Value = "  +1.39%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "1.002"
$ws.Range("E7").Value = "  -0.22%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.07404"
$ws.Range("E8").Value = "  -0.01%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.2946"
$ws.Range("E9").Value = "  -0.28%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "23.36"
$ws.Range("E10").Value = "  +1.46%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07678"
$ws.Range("E11").Value = "  -0.01%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.843.78"
$ws.Range("E12").Value = "  +0.95%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "5.026"
$ws.Range("E13").Value = "  +0.44%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.6767"
$ws.Range("E14").Value = "  +0.76%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "83.28"
$ws.Range("E15").Value = "  +0.54%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.000009171"
$ws.Range("E16").Value = "  +1.64%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "5.899"
$ws.Range("E17").Value = "  +0.09%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "29.496.74"
$ws.Range("E18").Value = "  +1.45%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "2.096.04"
$ws.Range("E19").Value = "  +1.23%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "242.93"
$ws.Range("E20").Value = "  +2.97%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "12.55"
$ws.Range("E21").Value = "  -1.04%  "
$ws.Range("E22").Value = "  -0.33%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "7.416"
$ws.Range("E23").Value = "  +3.18%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "1.002"
$ws.Range("E24").Value = "  -0.28%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "158.50"
$ws.Range("E25").Value = "  -0.53%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.1410"
$ws.Range("E26").Value = "  -1.76%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "8.525"
$ws.Range("E27").Value = "  +0.22%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "17.78"
$ws.Range("E28").Value = "  -0.31%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0.06152"
$ws.Range("E29").Value = "  +10.37%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.497"
$ws.Range("E30").Value = "  -0.02%  "
$ws.Range("B31").Value = "Toncoin"
$ws.Range("C31").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.232"
$ws.Range("E31").Value = "  +1.31%  "
$ws.Range("B32").Value = "Filecoin"
$ws.Range("C32").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "4.122"
$ws.Range("E32").Value = "  -0.32%  "
$ws.Range("B33").Value = "InternetComputer(DFINITY)"
$ws.Range("C33").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "4.101"
$ws.Range("E33").Value = "  -0.15%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.867"
$ws.Range("E34").Value = "  +1.01%  "
$ws.Range("B35").Value = "ARBITRUM"
$ws.Range("C35").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.144"
$ws.Range("E35").Value = "  +0.41%  "
$ws.Range("B36").Value = "ImmutableX"
$ws.Range("C36").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.7269"
$ws.Range("E36").Value = "  -2.46%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.613"
$ws.Range("E37").Value = "  -1.56%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.893"
$ws.Range("E38").Value = "  +3.83%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "1.225.80"
$ws.Range("E39").Value = "  +1.37%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.01768"
$ws.Range("E40").Value = "  -0.74%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "6.322"
$ws.Range("E41").Value = "  -0.18%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.9154"
$ws.Range("E42").Value = "  +2.23%  "
$ws.Range("E43").Value = "  -0.08%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "2.011.50"
$ws.Range("E44").Value = "  +1.76%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "101.82"
$ws.Range("E45").Value = "  +0.62%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "65.58"
$ws.Range("E46").Value = "  +0.65%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.00000000120"
$ws.Range("E47").Value = "  -0.58%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.5079"
$ws.Range("E48").Value = "  -0.42%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "9.256"
$ws.Range("E49").Value = "  +0.80%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.4062"
$ws.Range("E50").Value = "  +0.17%  "
$ws.Range("E51").Value = "  +4.78%  "
